# Wrote a pair of new test cases -> update counts, note, and selection on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two more partially automated test cases were added.
$ws.Range("E3").Value = "Contains four partially automated test cases."

# Total test cases and automated test cases both grew from 10 to 12.
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 12

# Move the active selection to reflect where the author ended up working.
$ws.Range("E4").Select()

$wb.Save()
